$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Cd200"
$ws.Range("C2").Value = "Cd200r1"
$ws.Range("D2").Value = "M1"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 70.0507755
$ws.Range("H2").Value = 140.101551
$ws.Range("I2").Value = 0.485408695871283
$ws.Range("J2").Value = 0.4148906539949664
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 15.53304966666667
$ws.Range("N2").Value = 46.599149
$ws.Range("O2").Value = 0.5918754867343928
$ws.Range("P2").Value = 0.5918754867343928
$ws.Range("Q2").Value = 1088.102175030016
$ws.Range("R2").Value = 6528.613050180099
$ws.Range("S2").Value = 0.2873015081339225
$ws.Range("T2").Value = 0.2455636077748213

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Cd200"
$ws.Range("C3").Value = "Cd200r1"
$ws.Range("D3").Value = "M2"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 70.0507755
$ws.Range("H3").Value = 140.101551
$ws.Range("I3").Value = 0.485408695871283
$ws.Range("J3").Value = 0.4148906539949664
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 10.71072966666667
$ws.Range("N3").Value = 32.132189
$ws.Range("O3").Value = 0.4081245132656072
$ws.Range("P3").Value = 0.4081245132656072
$ws.Range("Q3").Value = 750.2949193208566
$ws.Range("R3").Value = 4501.76951592514
$ws.Range("S3").Value = 0.1981071877373605
$ws.Range("T3").Value = 0.1693270462201451

# Row 4
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Cd200"
$ws.Range("C4").Value = "Cd200r1"
$ws.Range("D4").Value = "M1"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.7450486666666666
$ws.Range("H4").Value = 2.235146
$ws.Range("I4").Value = 0.005162728022151658
$ws.Range("J4").Value = 0.006619064379338907
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 15.53304966666667
$ws.Range("N4").Value = 46.599149
$ws.Range("O4").Value = 0.5918754867343928
$ws.Range("P4").Value = 0.5918754867343928
$ws.Range("Q4").Value = 11.57287794341711
$ws.Range("R4").Value = 104.155901490754
$ws.Range("S4").Value = 0.003055692160988302
$ws.Range("T4").Value = 0.003917661951247497

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Cd200"
$ws.Range("C5").Value = "Cd200r1"
$ws.Range("D5").Value = "M2"
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.7450486666666666
$ws.Range("H5").Value = 2.235146
$ws.Range("I5").Value = 0.005162728022151658
$ws.Range("J5").Value = 0.006619064379338907
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 10.71072966666667
$ws.Range("N5").Value = 32.132189
$ws.Range("O5").Value = 0.4081245132656072
$ws.Range("P5").Value = 0.4081245132656072
$ws.Range("Q5").Value = 7.98001485717711
$ws.Range("R5").Value = 71.82013371459401
$ws.Range("S5").Value = 0.002107035861163356
$ws.Range("T5").Value = 0.00270140242809141

# Row 6
$ws.Range("A6").Value = "M1"
$ws.Range("B6").Value = "Cd200"
$ws.Range("C6").Value = "Cd200r1"
$ws.Range("D6").Value = "M1"
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 0.6666666666666666
$ws.Range("G6").Value = 0.871287
$ws.Range("H6").Value = 2.613861
$ws.Range("I6").Value = 0.006037481860562736
$ws.Range("J6").Value = 0.007740574547543281
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 15.53304966666667
$ws.Range("N6").Value = 46.599149
$ws.Range("O6").Value = 0.5918754867343928
$ws.Range("P6").Value = 0.5918754867343928
$ws.Range("Q6").Value = 13.533744244921
$ws.Range("R6").Value = 121.803698204289
$ws.Range("S6").Value = 0.003573437514870637
$ws.Range("T6").Value = 0.004581456327931032

# Row 7
$ws.Range("A7").Value = "M1"
$ws.Range("B7").Value = "Cd200"
$ws.Range("C7").Value = "Cd200r1"
$ws.Range("D7").Value = "M2"
$ws.Range("E7").Value = 2
$ws.Range("F7").Value = 0.6666666666666666
$ws.Range("G7").Value = 0.871287
$ws.Range("H7").Value = 2.613861
$ws.Range("I7").Value = 0.006037481860562736
$ws.Range("J7").Value = 0.007740574547543281
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 10.71072966666667
$ws.Range("N7").Value = 32.132189
$ws.Range("O7").Value = 0.4081245132656072
$ws.Range("P7").Value = 0.4081245132656072
$ws.Range("Q7").Value = 9.332119519081001
$ws.Range("R7").Value = 83.98907567172901
$ws.Range("S7").Value = 0.002464044345692099
$ws.Range("T7").Value = 0.003159118219612249

# Row 8
$ws.Range("A8").Value = "M2"
$ws.Range("B8").Value = "Cd200"
$ws.Range("C8").Value = "Cd200r1"
$ws.Range("D8").Value = "M1"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 1.056899666666667
$ws.Range("H8").Value = 3.170699
$ws.Range("I8").Value = 0.007323663231443602
$ws.Range("J8").Value = 0.009389570439025233
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 15.53304966666667
$ws.Range("N8").Value = 46.599149
$ws.Range("O8").Value = 0.5918754867343928
$ws.Range("P8").Value = 0.5918754867343928
$ws.Range("Q8").Value = 16.41687501501677
$ws.Range("R8").Value = 147.751875135151
$ws.Range("S8").Value = 0.004334696739789459
$ws.Range("T8").Value = 0.005557456573824926

# Row 9
$ws.Range("A9").Value = "M2"
$ws.Range("B9").Value = "Cd200"
$ws.Range("C9").Value = "Cd200r1"
$ws.Range("D9").Value = "M2"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 1.056899666666667
$ws.Range("H9").Value = 3.170699
$ws.Range("I9").Value = 0.007323663231443602
$ws.Range("J9").Value = 0.009389570439025233
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 10.71072966666667
$ws.Range("N9").Value = 32.132189
$ws.Range("O9").Value = 0.4081245132656072
$ws.Range("P9").Value = 0.4081245132656072
$ws.Range("Q9").Value = 11.32016661445678
$ws.Range("R9").Value = 101.881499530111
$ws.Range("S9").Value = 0.002988966491654144
$ws.Range("T9").Value = 0.003832113865200307

# Row 10
$ws.Range("A10").Value = "Neutro"
$ws.Range("B10").Value = "Cd200"
$ws.Range("C10").Value = "Cd200r1"
$ws.Range("D10").Value = "M1"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 46.38387733333334
$ws.Range("H10").Value = 139.151632
$ws.Range("I10").Value = 0.321411679529899
$ws.Range("J10").Value = 0.4120776050862342
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 15.53304966666667
$ws.Range("N10").Value = 46.599149
$ws.Range("O10").Value = 0.5918754867343928
$ws.Range("P10").Value = 0.5918754867343928
$ws.Range("Q10").Value = 720.4830703512409
$ws.Range("R10").Value = 6484.347633161168
$ws.Range("S10").Value = 0.1902356942638776
$ws.Range("T10").Value = 0.2438986330827578

# Row 11
$ws.Range("A11").Value = "Neutro"
$ws.Range("B11").Value = "Cd200"
$ws.Range("C11").Value = "Cd200r1"
$ws.Range("D11").Value = "M2"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 46.38387733333334
$ws.Range("H11").Value = 139.151632
$ws.Range("I11").Value = 0.321411679529899
$ws.Range("J11").Value = 0.4120776050862342
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 10.71072966666667
$ws.Range("N11").Value = 32.132189
$ws.Range("O11").Value = 0.4081245132656072
$ws.Range("P11").Value = 0.4081245132656072
$ws.Range("Q11").Value = 496.805171009161
$ws.Range("R11").Value = 4471.246539082449
$ws.Range("S11").Value = 0.1311759852660213
$ws.Range("T11").Value = 0.1681789720034765

# Row 12
$ws.Range("A12").Value = "sCs"
$ws.Range("B12").Value = "Cd200"
$ws.Range("C12").Value = "Cd200r1"
$ws.Range("D12").Value = "M1"
$ws.Range("E12").Value = 2
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 25.205092
$ws.Range("H12").Value = 50.410184
$ws.Range("I12").Value = 0.17465575148466
$ws.Range("J12").Value = 0.1492825315528919
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 15.53304966666667
$ws.Range("N12").Value = 46.599149
$ws.Range("O12").Value = 0.5918754867343928
$ws.Range("P12").Value = 0.5918754867343928
$ws.Range("Q12").Value = 391.5119458889026
$ws.Range("R12").Value = 2349.071675333416
$ws.Range("S12").Value = 0.1033744579209443
$ws.Range("T12").Value = 0.08835667102381024

# Row 13
$ws.Range("A13").Value = "sCs"
$ws.Range("B13").Value = "Cd200"
$ws.Range("C13").Value = "Cd200r1"
$ws.Range("D13").Value = "M2"
$ws.Range("E13").Value = 2
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 25.205092
$ws.Range("H13").Value = 50.410184
$ws.Range("I13").Value = 0.17465575148466
$ws.Range("J13").Value = 0.1492825315528919
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 10.71072966666667
$ws.Range("N13").Value = 32.132189
$ws.Range("O13").Value = 0.4081245132656072
$ws.Range("P13").Value = 0.4081245132656072
$ws.Range("Q13").Value = 269.9649266354627
$ws.Range("R13").Value = 1619.789559812776
$ws.Range("S13").Value = 0.0712812935637157
$ws.Range("T13").Value = 0.06092586052908166
